# B6-PowerPoint.pptx edit — 2020-07-29
#
# 1) Three tables (slides 14, 15, 16) had their table style switched from
#    the deck's local "no style / grid" style to the built-in themed table
#    style {A8978370-055D-44C1-97AF-51E237D99A96}.
# 2) The presentation's main theme (color scheme carried on the slide
#    master's theme) was switched from the custom "Integral / Red Violet"
#    palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -------------------------------------
$newTableStyleId = "{A8978370-055D-44C1-97AF-51E237D99A96}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyleId)
    }
}

# --- 2. Swap the theme colour scheme to the standard Office palette ---
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Order matches the OOXML <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,          # dk1       000000
    16777215,   # lt1       FFFFFF
    6968388,    # dk2       44546A
    15132391,   # lt2       E7E6E6
    13998939,   # accent1   5B9BD5
    3243501,    # accent2   ED7D31
    10855845,   # accent3   A5A5A5
    49407,      # accent4   FFC000
    12874308,   # accent5   4472C4
    4697456,    # accent6   70AD47
    12673797,   # hlink     0563C1
    7491477     # folHlink  954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
